$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("H7").Value = 450.35
$ws.Range("H8").Select()
